# PFAS_Rdata_BW.xlsx - "Add files via upload" re-save.
#
# The authored change replaces the relative day labels in column A
# (D-2, D-1, D0, D1, D7 - stored as shared strings) with the actual
# calendar dates they represent, formatted as dates. The active
# selection also moved from F10 to E12 by the time the file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: turn the day-offset labels into real dates, formatted as
# short dates (maps to the built-in m/d/yyyy number format, id 14).
$ws.Range("A2:A6").NumberFormat = "mm-dd-yy"

$ws.Range("A2").Value = (Get-Date -Year 2025 -Month 3 -Day 20 -Hour 0 -Minute 0 -Second 0)
$ws.Range("A3").Value = (Get-Date -Year 2025 -Month 3 -Day 21 -Hour 0 -Minute 0 -Second 0)
$ws.Range("A4").Value = (Get-Date -Year 2025 -Month 3 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Range("A5").Value = (Get-Date -Year 2025 -Month 3 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("A6").Value = (Get-Date -Year 2025 -Month 4 -Day 1 -Hour 0 -Minute 0 -Second 0)

# Selection at save time moved to E12.
$ws.Range("E12").Select() | Out-Null
